$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Top WAR" ----
$ws1 = $wb.Worksheets.Item(1)

# Rename header B1 "Name" -> "Player"
$ws1.Range("B1").Value = "Player"

# Add new header C1 "WAR/pos", matching the existing header formatting
$ws1.Range("B1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)
$ws1.Range("C1").Value = "WAR/pos"

# Fill WAR/pos values for rows 2-11
$war = @(5.8, 5.6, 6.6, 3.7, 4, 3.4, 3.9, 6.4, 5.2, 7.3)
for ($i = 0; $i -lt $war.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 3).Value = $war[$i]
}

# ---- Sheet 2: "Top OPS Players" ----
$ws2 = $wb.Worksheets.Item(2)

# Add new headers D1 "OPS", E1 "AB", matching the existing header formatting
$ws2.Range("C1").Copy()
$ws2.Range("D1:E1").PasteSpecial(-4122)
$ws2.Range("D1").Value = "OPS"
$ws2.Range("E1").Value = "AB"

# Fill OPS (D) and AB (E) values for rows 2-11
$ops = @(0.9360000000000001, 0.915, 0.882, 0.867, 0.859, 0.851, 0.837, 0.836, 0.831, 0.827)
$ab  = @(590, 589, 655, 431, 597, 639, 635, 595, 573, 678)
for ($i = 0; $i -lt $ops.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 4).Value = $ops[$i]
    $ws2.Cells.Item($row, 5).Value = $ab[$i]
}
